$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1, Col 1: "165×8=1320" -> "142×3=426"
$t.Cell(1, 1).Range.Text = "142×3=426"
# Row 1, Col 2: "331×6=1986" -> "866×9=7794"
$t.Cell(1, 2).Range.Text = "866×9=7794"
# Row 1, Col 3: "823×5=4115" -> "603×5=3015"
$t.Cell(1, 3).Range.Text = "603×5=3015"
# Row 1, Col 4: "913×6=5478" -> "698×2=1396"
$t.Cell(1, 4).Range.Text = "698×2=1396"
# Row 1, Col 5: "971×8=7768" -> "881×2=1762"
$t.Cell(1, 5).Range.Text = "881×2=1762"
# Row 5, Col 1: "834×3=2502" -> "133×5=665"
$t.Cell(5, 1).Range.Text = "133×5=665"
# Row 5, Col 2: "603×8=4824" -> "129×2=258"
$t.Cell(5, 2).Range.Text = "129×2=258"
# Row 5, Col 3: "583×9=5247" -> "146×7=1022"
$t.Cell(5, 3).Range.Text = "146×7=1022"
# Row 5, Col 4: "622×6=3732" -> "722×2=1444"
$t.Cell(5, 4).Range.Text = "722×2=1444"
# Row 5, Col 5: "649×4=2596" -> "587×6=3522"
$t.Cell(5, 5).Range.Text = "587×6=3522"
# Row 10, Col 1: "586×2=1172" -> "953×4=3812"
$t.Cell(10, 1).Range.Text = "953×4=3812"
# Row 10, Col 2: "789×5=3945" -> "155×9=1395"
$t.Cell(10, 2).Range.Text = "155×9=1395"
# Row 10, Col 3: "477×9=4293" -> "987×9=8883"
$t.Cell(10, 3).Range.Text = "987×9=8883"
# Row 10, Col 4: "453×8=3624" -> "954×3=2862"
$t.Cell(10, 4).Range.Text = "954×3=2862"
# Row 10, Col 5: "677×8=5416" -> "145×2=290"
$t.Cell(10, 5).Range.Text = "145×2=290"
# Row 15, Col 1: "236×3=708" -> "708×9=6372"
$t.Cell(15, 1).Range.Text = "708×9=6372"
# Row 15, Col 2: "657×8=5256" -> "624×4=2496"
$t.Cell(15, 2).Range.Text = "624×4=2496"
# Row 15, Col 3: "678×3=2034" -> "350×2=700"
$t.Cell(15, 3).Range.Text = "350×2=700"
# Row 15, Col 4: "905×4=3620" -> "388×5=1940"
$t.Cell(15, 4).Range.Text = "388×5=1940"
# Row 15, Col 5: "816×2=1632" -> "778×4=3112"
$t.Cell(15, 5).Range.Text = "778×4=3112"
# Row 20, Col 1: "468×7=3276" -> "986×3=2958"
$t.Cell(20, 1).Range.Text = "986×3=2958"
# Row 20, Col 2: "708×9=6372" -> "550×4=2200"
$t.Cell(20, 2).Range.Text = "550×4=2200"
# Row 20, Col 3: "587×9=5283" -> "465×6=2790"
$t.Cell(20, 3).Range.Text = "465×6=2790"
# Row 20, Col 4: "377×3=1131" -> "401×2=802"
$t.Cell(20, 4).Range.Text = "401×2=802"
# Row 20, Col 5: "307×9=2763" -> "908×6=5448"
$t.Cell(20, 5).Range.Text = "908×6=5448"

Write-Host "Done applying cell updates"
